# C5-PowerPoint.pptx — commit: Sat, Jul 04, 2020 9:10:19 PM
#
# The only content-level change reachable through the PowerPoint object
# model is the table style applied to the "Sources of finance" table on
# slide 6: its <a:tableStyleId> moves from the custom "Table_0" style
# ({2B904C1B-C205-4DAA-AD67-B9F970814B58}, defined in tableStyles.xml) to
# the built-in style {9A2C342E-2271-4E0E-B399-AA17850B9077}.

$p = $ppt.ActivePresentation

$targetSlideIndex = 6
$newStyleId = "{9A2C342E-2271-4E0E-B399-AA17850B9077}"

$slide = $p.Slides.Item($targetSlideIndex)

# Find the shape that carries the table on this slide (it is shape 2, but
# look it up defensively in case shape ordering ever differs).
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shape = $slide.Shapes.Item($i)
    if ($shape.HasTable) {
        $shape.Table.ApplyStyle($newStyleId)
    }
}
